$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the id column (A) which was previously blank
$ws.Range("A2").Value = 1001
$ws.Range("A3").Value = 1002
$ws.Range("A4").Value = 1003

# New column K: product_discount
$ws.Range("K1").Value = "product_discount"
$ws.Range("K2").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("K4").Value = 0

# Apply the "default" (smaller) font style to H1, J1 and the new K column header/values
$ws.Range("H1").Font.Size = 11
$ws.Range("J1").Font.Size = 11
$ws.Range("K1:K4").Font.Size = 11

Write-Host "Done"
